$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: I1 = "I0", J1 = "IF", matching H1's style (bold/centered/bordered) ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows: columns I (I0) and J (IF) for rows 2-32 ---
$rows = @(
    @(2, 4, 8),
    @(3, 1, 5),
    @(4, 1, 6),
    @(5, 1, 5),
    @(6, 1, 6),
    @(7, 1, 4),
    @(8, 1, 6),
    @(9, 1, 7),
    @(10, 1, 7),
    @(11, 1, 7),
    @(12, 1, 5),
    @(13, 1, 6),
    @(14, 1, 7),
    @(15, 1, 7),
    @(16, 1, 6),
    @(17, 1, 7),
    @(18, 1, 7),
    @(19, 1, 5),
    @(20, 1, 6),
    @(21, 1, 6),
    @(22, 1, 4),
    @(23, 1, 6),
    @(24, 1, 6),
    @(25, 1, 5),
    @(26, 1, 5),
    @(27, 1, 6),
    @(28, 1, 5),
    @(29, 1, 5),
    @(30, 5, 8),
    @(31, 1, 3),
    @(32, 3, 4)
)

foreach ($item in $rows) {
    $r = $item[0]
    $i0 = $item[1]
    $if = $item[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}
